$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 178, shifting existing rows 178:237 down to 179:238
$ws.Rows("178:178").Insert()

# Populate the newly inserted row 178 with the new data record
$ws.Range("A178").Value = 4
$ws.Range("B178").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C178").Value = "Los Lagos"
$ws.Range("D178").Value = 44876
$ws.Range("E178").Value = 10
$ws.Range("F178").Value = 100112009
$ws.Range("G178").Value = "Acelga"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 80
$ws.Range("K178").Value = 10000
$ws.Range("L178").Value = 10000
$ws.Range("M178").Value = 10000
$ws.Range("N178").Value = "$/docena de atados (12 kilos)"
$ws.Range("O178").Value = "Región de La Araucanía"
$ws.Range("P178").Value = 833
$ws.Range("Q178").Value = 12
$ws.Range("R178").Value = "Hortaliza"
